$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    "59-2=57",
    "72+18=90",
    "57-49=8",
    "73-25=48",
    "47-30=17",
    "64-21=43",
    "62-53=9",
    "73+9=82",
    "98-15=83",
    "11+25=36",
    "73+18=91",
    "39-38=1",
    "59-3=56",
    "74-56=18",
    "9+24=33",
    "71-10=61",
    "21+11=32",
    "31+16=47",
    "51-31=20",
    "8+78=86",
    "7+61=68",
    "52-4=48",
    "33-18=15",
    "39+20=59",
    "74-69=5",
    "54+13=67",
    "21+68=89",
    "35+44=79",
    "60-13=47",
    "50-42=8",
    "9+29=38",
    "13+8=21",
    "37-22=15",
    "54+5=59",
    "92-36=56",
    "78+3=81",
    "77-3=74",
    "17+8=25",
    "40-23=17",
    "89-63=26",
    "9+26=35",
    "68+18=86",
    "65+11=76",
    "31+41=72",
    "27+69=96",
    "26+5=31",
    "18-11=7",
    "6+71=77",
    "41-22=19",
    "86+8=94",
    "68-0=68",
    "50-3=47",
    "6+1=7",
    "26+16=42",
    "32+42=74",
    "11-3=8",
    "30+33=63",
    "8+64=72",
    "11+5=16",
    "14+27=41",
    "80+7=87",
    "39-37=2",
    "14+61=75",
    "30-18=12",
    "45-26=19",
    "39-28=11",
    "80-3=77",
    "41-31=10",
    "3+20=23",
    "32+39=71",
    "28+4=32",
    "47+24=71",
    "75-61=14",
    "40+49=89",
    "91-26=65",
    "2+90=92",
    "98-80=18",
    "11+44=55",
    "92-24=68",
    "65-40=25",
    "1+94=95",
    "87+7=94",
    "38-30=8",
    "97-86=11",
    "92-37=55",
    "21+22=43",
    "28+15=43",
    "99-51=48",
    "40-37=3",
    "18-15=3",
    "54+8=62",
    "28-20=8",
    "68-59=9",
    "36+6=42",
    "65+3=68",
    "68+10=78",
    "17+56=73",
    "14+21=35",
    "29+26=55",
    "56+11=67"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$i]
        $i = $i + 1
    }
}

Write-Output ("Updated " + $i + " cells")
